$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("D11").Value = -7.798999999999998
$ws.Range("B12").Value = 5.041799999999998
$ws.Range("D23").Value = -8.291900000000002
$ws.Range("B27").Value = 6.202300000000004
$ws.Range("D28").Value = -8.096799999999995
$ws.Range("B32").Value = 6.195999999999999
$ws.Range("D32").Value = -8.380499999999985
$ws.Range("D34").Value = -8.0648
$ws.Range("B36").Value = 9.100700000000003
$ws.Range("B38").Value = 5.105999999999997
$ws.Range("D42").Value = -9.095199999999993
$ws.Range("B46").Value = 6.160900000000005
$ws.Range("D49").Value = -8.133599999999999
$ws.Range("B54").Value = 4.986600000000006
$ws.Range("D54").Value = -7.951199999999993
$ws.Range("B55").Value = 6.129799999999998
$ws.Range("B56").Value = 4.618899999999996
$ws.Range("B67").Value = 5.394499999999997
$ws.Range("B69").Value = 5.380899999999996
$ws.Range("B72").Value = 5.157600000000003
$ws.Range("D78").Value = -7.800800000000001
$ws.Range("D80").Value = -7.698099999999999
$ws.Range("B83").Value = 5.285199999999996
$ws.Range("B86").Value = 5.261900000000005
$ws.Range("B91").Value = 5.261899999999994
$ws.Range("B93").Value = 5.4261
$ws.Range("D97").Value = -8.30089999999999
$ws.Range("B99").Value = 4.830199999999999
$ws.Range("D99").Value = -8.404700000000002
$ws.Range("D101").Value = -7.825199999999999
$ws.Range("B104").Value = 10.029
